$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CodeModule" column (previously column F, with values AP61..AP69) moves to
# become column A, pushing ModuleName/ElementName1/ElementName2/Dept_Attachement/
# Coordinator one column to the right (A->B, B->C, C->D, D->E, E->F). At the same
# time the code values themselves are renamed from AP6x to G3EI3x.
$rows = @(
  @('CodeModule','ModuleName','ElementName1','ElementName2','Dept_Attachement','Coordinator'),
  @('G3EI31','Automatisme','Automatismes industriels','Supervision','SIC','Chater'),
  @('G3EI32','Technologie','RDM','conception','GEI','FILALI'),
  @('G3EI33','Solaire Thermique','Solaire thermique','Geothermie','GEI','EL FADAR'),
  @('G3EI34','Energie Eolienne','Aerodynamique','Etude economique','SIC','motaki'),
  @('G3EI35','Legislation ','Techniques EIE','projet etude','SIC','KAMACH'),
  @('G3EI36','Communication2','Allemand','Anglais','SIC','Haris'),
  @('G3EI37','Technologies','Technologie gaziere','Machines Thermiques','SIC','SARSRI'),
  @('G3EI38','Management','Projet ISO','SME','SIC','Khouya'),
  @('G3EI39','Maintenace et qualite','outils qualite','Maintenance surete','SIC','El kalkha')
)

$cols = @('A','B','C','D','E','F')

for ($r = 0; $r -lt $rows.Count; $r++) {
  $rowData = $rows[$r]
  for ($c = 0; $c -lt $cols.Count; $c++) {
    $addr = $cols[$c] + ($r + 1)
    $ws.Range($addr).Value = $rowData[$c]
  }
}

# Reflect the new active selection left behind in the saved sheet view.
$null = $ws.Range('B6').Select()
